# Insert a new weekly price record at row 364 for
# "Hortaliza, Terminal La Palmera de La Serena - Espinaca".
# This pushes the existing rows 364:486 down to 365:487 (each row keeps
# its own previously-held values), and the new row 364 is populated with
# the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 364, shifting rows
# 364:486 down to 365:487.
$ws.Rows.Item(364).Insert()

# Populate the newly inserted row 364 with the new record's data.
$ws.Cells.Item(364, 1).Value2 = 8
$ws.Cells.Item(364, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(364, 3).Value2 = "Coquimbo"
$ws.Cells.Item(364, 4).Value2 = 45215
$ws.Cells.Item(364, 5).Value2 = 4
$ws.Cells.Item(364, 6).Value2 = 100112012
$ws.Cells.Item(364, 7).Value2 = "Espinaca"
$ws.Cells.Item(364, 8).Value2 = "Sin especificar"
$ws.Cells.Item(364, 9).Value2 = "Primera"
$ws.Cells.Item(364, 10).Value2 = 1400
$ws.Cells.Item(364, 11).Value2 = 450
$ws.Cells.Item(364, 12).Value2 = 500
$ws.Cells.Item(364, 13).Value2 = 475
$ws.Cells.Item(364, 14).Value2 = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(364, 15).Value2 = "Provincia del Elqu" + [char]0x00ED
$ws.Cells.Item(364, 16).Value2 = 950
$ws.Cells.Item(364, 17).Value2 = 0.5
$ws.Cells.Item(364, 18).Value2 = "Hortaliza"

# Make sure the date column keeps the same date/time number format as the
# rest of column D.
$ws.Cells.Item(364, 4).NumberFormat = $ws.Cells.Item(365, 4).NumberFormat
